# Update cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text, matching the
# original inline-string cell type used throughout column D.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.394.38'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.876.14'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '0.7171'
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").Value = '239.84'
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.07830'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = '0.3092'
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("D10").Value = '24.75'
$ws.Range("E10").Value = '  +5.66%  '
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").Value = '1.869.02'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").Value = '0.7252'
$ws.Range("E13").Value = '  +3.35%  '
$ws.Range("D14").Value = '5.275'
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").Value = '91.14'
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = '29.415.21'
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '5.901'
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("E18").Value = '  +3.35%  '
$ws.Range("D19").Value = '0.000007897'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").Value = '13.28'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '7.916'
$ws.Range("E22").Value = '  +7.36%  '
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '0.1560'
$ws.Range("E24").Value = '  +8.97%  '
$ws.Range("D25").Value = '163.75'
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("D26").Value = '9.008'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  +1.53%  '
$ws.Range("E28").Value = '  -4.89%  '
$ws.Range("D29").Value = '1.488'
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").Value = '4.377'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").Value = '4.123'
$ws.Range("E31").Value = '  +2.80%  '
$ws.Range("D32").Value = '0.05279'
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("D33").Value = '1.926'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '1.199'
$ws.Range("E34").Value = '  +3.21%  '
$ws.Range("D35").Value = '0.7209'
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("D36").Value = '2.676'
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("D37").Value = '0.01860'
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("D38").Value = '1.224.42'
$ws.Range("E38").Value = '  +9.14%  '
$ws.Range("D39").Value = '2.718'
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = '0.9090'
$ws.Range("E40").Value = '  -2.20%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '72.72'
$ws.Range("E41").Value = '  +4.45%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.069'
$ws.Range("E42").Value = '  +4.13%  '
$ws.Range("D43").Value = '0.9999'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '103.60'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = '0.5340'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").Value = '0.00000000121'
$ws.Range("E46").Value = '  +3.99%  '
$ws.Range("D47").Value = '2.937'
$ws.Range("E47").Value = '  +11.75%  '
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").Value = '0.4319'
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.07%  '
